# "Generate Report for Archive"
# - Flip the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F2, and the
#   "Status" column on each per-locale sheet).
# - The shorter replacement text lets Excel's column auto-sizing narrow
#   the affected "Status" columns, so nudge their widths down to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the now-shorter "Status" columns (Overview columns E/F, and
# column C on each locale sheet) to their new auto-fit width.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
